$wb = $excel.ActiveWorkbook

# Update User ID & Password on the Login sheet
$login = $wb.Worksheets.Item("Login")
$login.Range("A2").Value = "mngr516795"
$login.Range("B2").Value = "hYtYvYz"

# Make Login the active sheet and select A5 (mirrors the commit's view-state change)
$login.Activate()
$login.Range("A5").Select()
